$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column values
$ws.Range("H4").Value = 2
$ws.Range("H7").Value = 6
$ws.Range("H8").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H18").Value = 7
$ws.Range("H19").Value = 5
$ws.Range("H20").Value = 5

# Update view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("I20").Select()
